$wb = $excel.ActiveWorkbook

# This script applies the numeric corrections described in the commit diff
# ("chore: update Sheets via scheduled runner") to the underlying leve-profit
# tables on each sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1346.7
$ws.Range("I6").Value = 1346.7
$ws.Range("K6").Value = 4040.1
$ws.Range("M6").Value = -3928.1
$ws.Range("H29").Value = 100
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H112").Value = 580428.5600000001
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 580428.5600000001
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 1741285.68
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -1743501.68
$ws.Range("H121").Value = 1124.1666
$ws.Range("I121").Value = 575
$ws.Range("J121").Value = 1398.75
$ws.Range("K121").Value = 1725
$ws.Range("L121").Value = 4196.25
$ws.Range("M121").Value = 22
$ws.Range("N121").Value = -7690.25
$ws.Range("H137").Value = 14893702
$ws.Range("I137").Value = 1046.9667
$ws.Range("J137").Value = 52125340
$ws.Range("K137").Value = 3140.9001
$ws.Range("L137").Value = 156376020
$ws.Range("M137").Value = -590.9000999999998
$ws.Range("N137").Value = -156381120
$ws.Range("H138").Value = 2167.1462
$ws.Range("I138").Value = 1580.9788
$ws.Range("J138").Value = 2954.2856
$ws.Range("K138").Value = 4742.936400000001
$ws.Range("L138").Value = 8862.856800000001
$ws.Range("M138").Value = 397.0635999999995
$ws.Range("N138").Value = -19142.8568
$ws.Range("H141").Value = 1602.14
$ws.Range("I141").Value = 890.7442
$ws.Range("J141").Value = 5972.143
$ws.Range("K141").Value = 2672.2326
$ws.Range("L141").Value = 17916.429
$ws.Range("M141").Value = 2507.7674
$ws.Range("N141").Value = -28276.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 45099628
$ws.Range("I74").Value = 35715028
$ws.Range("K74").Value = 35715028
$ws.Range("M74").Value = -35714154
$ws.Range("H77").Value = 45099628
$ws.Range("I77").Value = 35715028
$ws.Range("K77").Value = 178575140
$ws.Range("M77").Value = -178570772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1881.46
$ws.Range("I86").Value = 1909.5667
$ws.Range("J86").Value = 1628.5
$ws.Range("K86").Value = 1909.5667
$ws.Range("L86").Value = 1628.5
$ws.Range("M86").Value = -786.5667000000001
$ws.Range("N86").Value = -3874.5
$ws.Range("H89").Value = 1881.46
$ws.Range("I89").Value = 1909.5667
$ws.Range("J89").Value = 1628.5
$ws.Range("K89").Value = 9547.833500000001
$ws.Range("L89").Value = 8142.5
$ws.Range("M89").Value = -3931.833500000001
$ws.Range("N89").Value = -19374.5
$ws.Range("H134").Value = 14881811
$ws.Range("I134").Value = 15625853
$ws.Range("K134").Value = 46877559
$ws.Range("M134").Value = -46875024

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 266.06668
$ws.Range("I5").Value = 86.375
$ws.Range("J5").Value = 471.42856
$ws.Range("K5").Value = 86.375
$ws.Range("L5").Value = 471.42856
$ws.Range("M5").Value = 25.625
$ws.Range("N5").Value = -695.4285600000001
$ws.Range("H31").Value = 1841864.5
$ws.Range("I31").Value = 1232.7084
$ws.Range("J31").Value = 6259381
$ws.Range("K31").Value = 1232.7084
$ws.Range("L31").Value = 6259381
$ws.Range("M31").Value = -937.7084
$ws.Range("N31").Value = -6259971
$ws.Range("H34").Value = 1841864.5
$ws.Range("I34").Value = 1232.7084
$ws.Range("J34").Value = 6259381
$ws.Range("K34").Value = 1232.7084
$ws.Range("L34").Value = 6259381
$ws.Range("M34").Value = -1030.7084
$ws.Range("N34").Value = -6259785
$ws.Range("H105").Value = 5206.3125
$ws.Range("I105").Value = 1475
$ws.Range("J105").Value = 16400.25
$ws.Range("K105").Value = 1475
$ws.Range("L105").Value = 16400.25
$ws.Range("M105").Value = 272
$ws.Range("N105").Value = -19894.25
$ws.Range("H132").Value = 1221.6
$ws.Range("I132").Value = 1166
$ws.Range("K132").Value = 3498
$ws.Range("M132").Value = -968

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 190788.06
$ws.Range("I7").Value = 277330.8
$ws.Range("J7").Value = 394
$ws.Range("K7").Value = 831992.3999999999
$ws.Range("L7").Value = 1182
$ws.Range("M7").Value = -831880.3999999999
$ws.Range("N7").Value = -1406
$ws.Range("H80").Value = 9475
$ws.Range("I80").Value = 950
$ws.Range("J80").Value = 11180
$ws.Range("K80").Value = 2850
$ws.Range("L80").Value = 33540
$ws.Range("M80").Value = -1914
$ws.Range("N80").Value = -35412
$ws.Range("H83").Value = 9475
$ws.Range("I83").Value = 950
$ws.Range("J83").Value = 11180
$ws.Range("K83").Value = 8550
$ws.Range("L83").Value = 100620
$ws.Range("M83").Value = -3870
$ws.Range("N83").Value = -109980
$ws.Range("H92").Value = 1666284
$ws.Range("I92").Value = 202
$ws.Range("J92").Value = 3570377.8
$ws.Range("K92").Value = 606
$ws.Range("L92").Value = 10711133.4
$ws.Range("M92").Value = 642
$ws.Range("N92").Value = -10713629.4
$ws.Range("H121").Value = 1594924.8
$ws.Range("I121").Value = 356
$ws.Range("J121").Value = 1776125.8
$ws.Range("K121").Value = 1068
$ws.Range("L121").Value = 5328377.4
$ws.Range("M121").Value = 242
$ws.Range("N121").Value = -5330997.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 2251.1667
$ws.Range("I9").Value = 1051.75
$ws.Range("J9").Value = 4650
$ws.Range("K9").Value = 1051.75
$ws.Range("L9").Value = 4650
$ws.Range("M9").Value = -881.75
$ws.Range("N9").Value = -4990
$ws.Range("H132").Value = 22569450
$ws.Range("I132").Value = 30953122
$ws.Range("J132").Value = 12988109
$ws.Range("K132").Value = 92859366
$ws.Range("L132").Value = 38964327
$ws.Range("M132").Value = -92856836
$ws.Range("N132").Value = -38969387

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2709.6667
$ws.Range("I7").Value = 2051.6
$ws.Range("K7").Value = 2051.6
$ws.Range("M7").Value = -1939.6
$ws.Range("H40").Value = 2410.2778
$ws.Range("I40").Value = 2491.1538
$ws.Range("J40").Value = 2200
$ws.Range("K40").Value = 2491.1538
$ws.Range("L40").Value = 2200
$ws.Range("M40").Value = -2355.1538
$ws.Range("N40").Value = -2472
$ws.Range("H126").Value = 2709.6667
$ws.Range("I126").Value = 2051.6
$ws.Range("K126").Value = 6154.799999999999
$ws.Range("M126").Value = -3684.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 5724
$ws.Range("J29").Value = 5724
$ws.Range("L29").Value = 5724
$ws.Range("N29").Value = -6304
$ws.Range("H49").Value = 11230
$ws.Range("J49").Value = 13980
$ws.Range("L49").Value = 13980
$ws.Range("N49").Value = -14440
$ws.Range("H113").Value = 388.85715
$ws.Range("I113").Value = 327.36365
$ws.Range("J113").Value = 614.3333
$ws.Range("K113").Value = 982.09095
$ws.Range("L113").Value = 1842.9999
$ws.Range("M113").Value = 1187.90905
$ws.Range("N113").Value = -6182.9999

